$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues constant (used to collapse a temporary formula down to a
# literal value/shared string without pulling in a new number format /
# cell style - plain .Value assignment would silently coerce numeric-
# looking text like "102" into a real number).
$xlPasteValues = -4163

function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = "=""$escaped"""
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

# Row 1 was the header ("789","4543","11231","2018-11-21"); its cells are
# now blanked out while the row itself (and its custom formatting) stays.
$ws.Range("A1:D1").ClearContents()

# Rows 2-4 are brand new rows carrying the same repeating record that the
# rest of the sheet already uses. They leave the old header's custom row
# formatting behind (plain data rows, like 5-10 below).
$newRows = 2, 3, 4
foreach ($r in $newRows) {
    $ws.Rows.Item($r).ClearFormats()
    $ws.Range("A$r").Value = "od-856429"
    Set-TextValue $ws.Range("B$r") "102"
    $ws.Range("C$r").Value = "user101"
    $ws.Range("D$r").Value = "2019-Sep-13"
}

# Row 5's order id now carries two leading spaces; row 8's loses them.
Set-TextValue $ws.Range("B5") "  102"
Set-TextValue $ws.Range("B8") "102"

# The old trailing row 11 is no longer needed - its data now lives in the
# newly added rows 2-4, so it is cleared out entirely.
$ws.Range("A11:D11").ClearContents()
